$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 62803.875
$ws.Range("I38").Value = 71658
$ws.Range("J38").Value = 825
$ws.Range("K38").Value = 214974
$ws.Range("L38").Value = 2475
$ws.Range("M38").Value = -214602
$ws.Range("N38").Value = -3219

$ws.Range("H39").Value = 14525.571
$ws.Range("I39").Value = 16866.5
$ws.Range("J39").Value = 480
$ws.Range("K39").Value = 50599.5
$ws.Range("L39").Value = 1440
$ws.Range("M39").Value = -50303.5
$ws.Range("N39").Value = -2032

$ws.Range("H40").Value = 1960
$ws.Range("I40").Value = 2400
$ws.Range("J40").Value = 1666.6666
$ws.Range("K40").Value = 2400
$ws.Range("L40").Value = 1666.6666
$ws.Range("M40").Value = -2225
$ws.Range("N40").Value = -2016.6666

$ws.Range("H43").Value = 946.3333
$ws.Range("I43").Value = 794.8
$ws.Range("K43").Value = 794.8
$ws.Range("M43").Value = -725.8

$ws.Range("H58").Value = 1770.7273
$ws.Range("I58").Value = 1258.5714
$ws.Range("J58").Value = 2667
$ws.Range("K58").Value = 3775.7142
$ws.Range("L58").Value = 8001
$ws.Range("M58").Value = -3625.7142
$ws.Range("N58").Value = -8301

$ws.Range("H100").Value = 1531.5358
$ws.Range("I100").Value = 1683.1666
$ws.Range("J100").Value = 1258.6
$ws.Range("K100").Value = 1683.1666
$ws.Range("L100").Value = 1258.6
$ws.Range("M100").Value = -1142.1666
$ws.Range("N100").Value = -2340.6

$ws.Range("H123").Value = 100780
$ws.Range("J123").Value = 100780
$ws.Range("L123").Value = 100780
$ws.Range("N123").Value = -110580

$ws.Range("H137").Value = 1475.05
$ws.Range("I137").Value = 1095.2727
$ws.Range("K137").Value = 3285.8181
$ws.Range("M137").Value = -735.8181

$ws.Range("H141").Value = 3515.9688
$ws.Range("I141").Value = 1891.3478
$ws.Range("J141").Value = 7667.778
$ws.Range("K141").Value = 5674.0434
$ws.Range("L141").Value = 23003.334
$ws.Range("M141").Value = -494.0434000000005
$ws.Range("N141").Value = -33363.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 362428.78
$ws.Range("I32").Value = 447323.44
$ws.Range("J32").Value = 12238.4375
$ws.Range("K32").Value = 447323.44
$ws.Range("L32").Value = 12238.4375
$ws.Range("M32").Value = -447036.44
$ws.Range("N32").Value = -12812.4375

$ws.Range("H43").Value = 9170.5
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H61").Value = 9806365
$ws.Range("I61").Value = 23811154
$ws.Range("J61").Value = 3012.9
$ws.Range("K61").Value = 23811154
$ws.Range("L61").Value = 3012.9
$ws.Range("M61").Value = -23810942
$ws.Range("N61").Value = -3436.9

$ws.Range("H74").Value = 732.19446
$ws.Range("I74").Value = 640.2308
$ws.Range("J74").Value = 784.1739
$ws.Range("K74").Value = 640.2308
$ws.Range("L74").Value = 784.1739
$ws.Range("M74").Value = 233.7692
$ws.Range("N74").Value = -2532.1739

$ws.Range("H77").Value = 732.19446
$ws.Range("I77").Value = 640.2308
$ws.Range("J77").Value = 784.1739
$ws.Range("K77").Value = 3201.154
$ws.Range("L77").Value = 3920.8695
$ws.Range("M77").Value = 1166.846
$ws.Range("N77").Value = -12656.8695

$ws.Range("H136").Value = 9806365
$ws.Range("I136").Value = 23811154
$ws.Range("J136").Value = 3012.9
$ws.Range("K136").Value = 71433462
$ws.Range("L136").Value = 9038.700000000001
$ws.Range("M136").Value = -71430912
$ws.Range("N136").Value = -14138.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1640.9032
$ws.Range("I20").Value = 1844.7646
$ws.Range("J20").Value = 1393.3572
$ws.Range("K20").Value = 1844.7646
$ws.Range("L20").Value = 1393.3572
$ws.Range("M20").Value = -1597.7646
$ws.Range("N20").Value = -1887.3572

$ws.Range("H64").Value = 1033.7142
$ws.Range("J64").Value = 1116.75
$ws.Range("L64").Value = 1116.75
$ws.Range("N64").Value = -1566.75

$ws.Range("H67").Value = 1033.7142
$ws.Range("J67").Value = 1116.75
$ws.Range("L67").Value = 1116.75
$ws.Range("N67").Value = -2676.75

$ws.Range("H134").Value = 2275.7114
$ws.Range("I134").Value = 2098.925
$ws.Range("K134").Value = 6296.775000000001
$ws.Range("M134").Value = -3761.775000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2228.8696
$ws.Range("I31").Value = 1851.6923
$ws.Range("J31").Value = 2719.2
$ws.Range("K31").Value = 1851.6923
$ws.Range("L31").Value = 2719.2
$ws.Range("M31").Value = -1556.6923
$ws.Range("N31").Value = -3309.2

$ws.Range("H34").Value = 2228.8696
$ws.Range("I34").Value = 1851.6923
$ws.Range("J34").Value = 2719.2
$ws.Range("K34").Value = 1851.6923
$ws.Range("L34").Value = 2719.2
$ws.Range("M34").Value = -1649.6923
$ws.Range("N34").Value = -3123.2

$ws.Range("H58").Value = 3599.111
$ws.Range("I58").Value = 3599.111
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3599.111
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -3396.111
$ws.Range("N58").ClearContents()

$ws.Range("H136").Value = 3599.111
$ws.Range("I136").Value = 3599.111
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10797.333
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8247.332999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 901.56665
$ws.Range("I68").Value = 945
$ws.Range("J68").Value = 900.5795000000001
$ws.Range("K68").Value = 2835
$ws.Range("L68").Value = 2701.7385
$ws.Range("M68").Value = -2024
$ws.Range("N68").Value = -4323.7385

$ws.Range("H71").Value = 901.56665
$ws.Range("I71").Value = 945
$ws.Range("J71").Value = 900.5795000000001
$ws.Range("K71").Value = 8505
$ws.Range("L71").Value = 8105.2155
$ws.Range("M71").Value = -4449
$ws.Range("N71").Value = -16217.2155

$ws.Range("H80").Value = 3375.0625
$ws.Range("I80").Value = 3000.25
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 9000.75
$ws.Range("L80").Value = 10500
$ws.Range("M80").Value = -8064.75
$ws.Range("N80").Value = -12372

$ws.Range("H83").Value = 3375.0625
$ws.Range("I83").Value = 3000.25
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 27002.25
$ws.Range("L83").Value = 31500
$ws.Range("M83").Value = -22322.25
$ws.Range("N83").Value = -40860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 912.75
$ws.Range("I107").Value = 475.25
$ws.Range("J107").Value = 1350.25
$ws.Range("K107").Value = 475.25
$ws.Range("L107").Value = 1350.25
$ws.Range("M107").Value = 1444.75
$ws.Range("N107").Value = -5190.25

$ws.Range("H122").Value = 1490
$ws.Range("I122").Value = 1490
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4470
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2020
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2636.4285
$ws.Range("I132").Value = 2279.8845
$ws.Range("J132").Value = 3666.4443
$ws.Range("K132").Value = 6839.6535
$ws.Range("L132").Value = 10999.3329
$ws.Range("M132").Value = -4309.6535
$ws.Range("N132").Value = -16059.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4334.5386
$ws.Range("I132").Value = 4186.5386
$ws.Range("J132").Value = 4482.5386
$ws.Range("K132").Value = 12559.6158
$ws.Range("L132").Value = 13447.6158
$ws.Range("M132").Value = -10029.6158
$ws.Range("N132").Value = -18507.6158

$ws.Range("H136").Value = 1300.0834
$ws.Range("I136").Value = 1257.2106
$ws.Range("K136").Value = 3771.6318
$ws.Range("M136").Value = -1221.6318

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3219.5527
$ws.Range("I136").Value = 3229.2273
$ws.Range("J136").Value = 3206.25
$ws.Range("K136").Value = 9687.6819
$ws.Range("L136").Value = 9618.75
$ws.Range("M136").Value = -7137.6819
$ws.Range("N136").Value = -14718.75
